$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.06
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 3.45
$ws.Range("V2").Value = 1.4
$ws.Range("X2").Value = 980
$ws.Range("Y2").Value = 26
$ws.Range("Z2").Value = 980
$ws.Range("AA2").Value = 65
$ws.Range("AB2").Value = 18.5
$ws.Range("AC2").Value = 13
$ws.Range("AD2").Value = 16.5
$ws.Range("AE2").Value = 1000
$ws.Range("AH2").Value = 16
$ws.Range("AJ2").Value = 32
$ws.Range("AK2").Value = 22
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 9.199999999999999
$ws.Range("AO2").Value = 24
$ws.Range("G3").Value = 2.2
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 3.95
$ws.Range("N3").Value = 4.8
$ws.Range("R3").Value = 1.52
$ws.Range("S3").Value = 2.32
$ws.Range("T3").Value = 1.58
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.33
$ws.Range("X3").Value = 27
$ws.Range("Y3").Value = 19
$ws.Range("Z3").Value = 980
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 13.5
$ws.Range("AC3").Value = 10
$ws.Range("AD3").Value = 16.5
$ws.Range("AE3").Value = 980
$ws.Range("AF3").Value = 16
$ws.Range("AG3").Value = 11.5
$ws.Range("AH3").Value = 16.5
$ws.Range("AI3").Value = 980
$ws.Range("AJ3").Value = 980
$ws.Range("AK3").Value = 21
$ws.Range("AL3").Value = 980
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 11.5
$ws.Range("G4").Value = 6.6
$ws.Range("H4").Value = 1.64
$ws.Range("I4").Value = 1.76
$ws.Range("N4").Value = 3.75
$ws.Range("P4").Value = 1.96
$ws.Range("V4").Value = 2.3
$ws.Range("W4").Value = 1.18
$ws.Range("X4").Value = 19
$ws.Range("Y4").Value = 8.800000000000001
$ws.Range("AA4").Value = 18
$ws.Range("AC4").Value = 9.6
$ws.Range("AE4").Value = 19
$ws.Range("AF4").Value = 1000
$ws.Range("AG4").Value = 24
$ws.Range("AI4").Value = 980
$ws.Range("AK4").Value = 85
$ws.Range("AL4").Value = 1000
$ws.Range("AM4").Value = 140
$ws.Range("AN4").Value = 1000
$ws.Range("N5").Value = 1.11
$ws.Range("S5").Value = 1.05
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 3.25
$ws.Range("H6").Value = 2.64
$ws.Range("I6").Value = 2.86
$ws.Range("N6").Value = 2.56
$ws.Range("P6").Value = 1.54
$ws.Range("R6").Value = 1.19
$ws.Range("S6").Value = 5.4
$ws.Range("T6").Value = 2.1
$ws.Range("U6").Value = 1.78
$ws.Range("V6").Value = 1.55
$ws.Range("X6").Value = 10.5
$ws.Range("AH6").Value = 1000
$ws.Range("AJ6").Value = 1000
$ws.Range("AN6").Value = 60
$ws.Range("AO6").Value = 1000
